$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
# Try inserting rows
$ws.Rows("108:110").Insert()
Write-Output "inserted"
